# Notenrechner Codereview SWE 2 - "Doku + kleine Verbesserungen + Add Photographer hinzugeüfgt"
#
# Updates the point values ("Punkte", column B) for several criteria rows
# on the "Notenrechner" sheet and adjusts the view (zoom / scroll / selection)
# to match where the reviewer ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notenrechner")
$ws.Activate()

# "Bilder zu Fotografen_innen zuordnen" bekam die volle Punktzahl (0.5 -> 1)
$ws.Range("B16").Value = 1

# Zuvor leere Bewertungen wurden explizit mit 0 Punkten befüllt
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0

# "Dokumentation vorhanden" bekam die volle Punktzahl
$ws.Range("B27").Value = 1

# View-Status: reingezoomt und near B33/D30 weitergearbeitet
$excel.ActiveWindow.Zoom = 117
$ws.Range("D30").Select()
